# Auto-generated Excel COM-interop script
# Applies numeric corrections to Kujata_Profits across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 302.4
$ws.Range("I38").Value = 302.4
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 907.1999999999999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -535.1999999999999

$ws.Range("H40").Value = 2771.4285
$ws.Range("I40").Value = 3180
$ws.Range("J40").Value = 1750
$ws.Range("K40").Value = 3180
$ws.Range("L40").Value = 1750
$ws.Range("M40").Value = -3005
$ws.Range("N40").Value = -2100

$ws.Range("H58").Value = 533.1667

$ws.Range("H69").Value = 2996
$ws.Range("I69").Value = 2713
$ws.Range("J69").Value = 3036.4285
$ws.Range("K69").Value = 8139
$ws.Range("L69").Value = 9109.2855
$ws.Range("M69").Value = -7265
$ws.Range("N69").Value = -10857.2855

$ws.Range("H72").Value = 2996
$ws.Range("I72").Value = 2713
$ws.Range("J72").Value = 3036.4285
$ws.Range("K72").Value = 24417
$ws.Range("L72").Value = 27327.8565
$ws.Range("M72").Value = -20049
$ws.Range("N72").Value = -36063.8565

$ws.Range("H74").Value = 4000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 4000
$ws.Range("N74").Value = -5872

$ws.Range("H77").Value = 4000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 20000
$ws.Range("N77").Value = -29360

$ws.Range("H92").Value = 719.36365
$ws.Range("I92").Value = 596.1053000000001
$ws.Range("J92").Value = 1500
$ws.Range("K92").Value = 596.1053000000001
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 651.8946999999999

$ws.Range("H112").Value = 1962.4681
$ws.Range("I112").Value = 705.9
$ws.Range("J112").Value = 2302.081
$ws.Range("K112").Value = 2117.7
$ws.Range("L112").Value = 6906.243
$ws.Range("M112").Value = -1009.7
$ws.Range("N112").Value = -9122.243

$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents() | Out-Null

$ws.Range("H133").Value = 29935.8
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 29935.8
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 29935.8
$ws.Range("N133").Value = -40055.8

$ws.Range("H138").Value = 519322.1
$ws.Range("I138").Value = 1437.4445
$ws.Range("J138").Value = 764635.9
$ws.Range("K138").Value = 4312.333500000001
$ws.Range("L138").Value = 2293907.7
$ws.Range("M138").Value = 827.6664999999994
$ws.Range("N138").Value = -2304187.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 850
$ws.Range("I25").Value = 850
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 850
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -448

$ws.Range("H28").Value = 4556
$ws.Range("I28").Value = 4556
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 4556
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -4364

$ws.Range("H32").Value = 3110.048
$ws.Range("I32").Value = 2795.808
$ws.Range("J32").Value = 5404
$ws.Range("K32").Value = 2795.808
$ws.Range("L32").Value = 5404
$ws.Range("M32").Value = -2508.808
$ws.Range("N32").Value = -5978

$ws.Range("H35").Value = 2281.5
$ws.Range("I35").Value = 2281.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2281.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1875.5

$ws.Range("H45").Value = 1130.7391
$ws.Range("I45").Value = 1086.0625
$ws.Range("J45").Value = 1232.8572
$ws.Range("K45").Value = 1086.0625
$ws.Range("L45").Value = 1232.8572
$ws.Range("M45").Value = -709.0625

$ws.Range("H61").Value = 1507.3334
$ws.Range("I61").Value = 1328.091
$ws.Range("J61").Value = 2296
$ws.Range("K61").Value = 1328.091
$ws.Range("L61").Value = 2296
$ws.Range("M61").Value = -1116.091
$ws.Range("N61").Value = -2720

$ws.Range("H74").Value = 1468.45
$ws.Range("I74").Value = 869.4
$ws.Range("J74").Value = 3265.6
$ws.Range("K74").Value = 869.4
$ws.Range("L74").Value = 3265.6
$ws.Range("M74").Value = 4.600000000000023
$ws.Range("N74").Value = -5013.6

$ws.Range("H76").Value = 25999.75
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 25999.75
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 25999.75
$ws.Range("N76").Value = -26675.75

$ws.Range("H77").Value = 1468.45
$ws.Range("I77").Value = 869.4
$ws.Range("J77").Value = 3265.6
$ws.Range("K77").Value = 4347
$ws.Range("L77").Value = 16328
$ws.Range("M77").Value = 21
$ws.Range("N77").Value = -25064

$ws.Range("H79").Value = 25999.75
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 25999.75
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 25999.75
$ws.Range("N79").Value = -28339.75

$ws.Range("H99").Value = 4556
$ws.Range("I99").Value = 4556
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4556
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1561

$ws.Range("H122").Value = 946.2273
$ws.Range("I122").Value = 943.6667
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2831.0001
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -381.0001000000002

$ws.Range("H132").Value = 2442.7334
$ws.Range("I132").Value = 2049.4614
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6148.3842
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3618.3842
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 1507.3334
$ws.Range("I136").Value = 1328.091
$ws.Range("J136").Value = 2296
$ws.Range("K136").Value = 3984.273
$ws.Range("L136").Value = 6888
$ws.Range("M136").Value = -1434.273
$ws.Range("N136").Value = -11988

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 55556668
$ws.Range("I99").Value = 71429544
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 71429544
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = -71428046
$ws.Range("N99").Value = -4596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1166.2106
$ws.Range("I31").Value = 1151.3214
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 1151.3214
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -856.3214
$ws.Range("N31").Value = -2590

$ws.Range("H34").Value = 1166.2106
$ws.Range("I34").Value = 1151.3214
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 1151.3214
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -949.3214
$ws.Range("N34").Value = -2404

$ws.Range("H114").Value = 24376
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 24376
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 24376
$ws.Range("N114").Value = -33054

$ws.Range("H132").Value = 6645.9585
$ws.Range("I132").Value = 7815
$ws.Range("J132").Value = 3806.8572
$ws.Range("K132").Value = 23445
$ws.Range("L132").Value = 11420.5716
$ws.Range("M132").Value = -20915

$ws.Range("H134").Value = 2646.3333
$ws.Range("I134").Value = 3021.5
$ws.Range("J134").Value = 1896
$ws.Range("K134").Value = 9064.5
$ws.Range("L134").Value = 5688
$ws.Range("M134").Value = -6529.5
$ws.Range("N134").Value = -10758

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 538.4
$ws.Range("I23").Value = 1100
$ws.Range("J23").Value = 398
$ws.Range("K23").Value = 3300
$ws.Range("L23").Value = 1194
$ws.Range("M23").Value = -3065
$ws.Range("N23").Value = -1664

$ws.Range("H130").Value = 1809.2858
$ws.Range("I130").Value = 1250
$ws.Range("J130").Value = 2033
$ws.Range("K130").Value = 3750
$ws.Range("L130").Value = 6099
$ws.Range("M130").Value = 1270

$ws.Range("H131").Value = 12989350
$ws.Range("I131").Value = 500000260
$ws.Range("J131").Value = 2392.2932
$ws.Range("K131").Value = 1500000780
$ws.Range("L131").Value = 7176.8796
$ws.Range("M131").Value = -1499995740
$ws.Range("N131").Value = -17256.8796

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6429412
$ws.Range("I12").Value = 6393750
$ws.Range("J12").Value = 7000000
$ws.Range("K12").Value = 6393750
$ws.Range("L12").Value = 7000000
$ws.Range("M12").Value = -6393610

$ws.Range("H26").Value = 16000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 16000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 16000
$ws.Range("M26").ClearContents() | Out-Null
$ws.Range("N26").Value = -16560

$ws.Range("H50").Value = 16000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 16000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 16000
$ws.Range("M50").ClearContents() | Out-Null
$ws.Range("N50").Value = -16996

$ws.Range("H132").Value = 2226.0244
$ws.Range("I132").Value = 1876.7826
$ws.Range("J132").Value = 2672.2778
$ws.Range("K132").Value = 5630.3478
$ws.Range("L132").Value = 8016.8334
$ws.Range("M132").Value = -3100.3478
$ws.Range("N132").Value = -13076.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 250.52
$ws.Range("I55").Value = 212.53847
$ws.Range("J55").Value = 291.66666
$ws.Range("K55").Value = 212.53847
$ws.Range("L55").Value = 291.66666
$ws.Range("M55").Value = -39.53846999999999
$ws.Range("N55").Value = -637.66666

$ws.Range("H61").Value = 2485.6667
$ws.Range("I61").Value = 2004
$ws.Range("J61").Value = 2582
$ws.Range("K61").Value = 2004
$ws.Range("L61").Value = 2582
$ws.Range("M61").Value = -1802
$ws.Range("N61").Value = -2986

$ws.Range("H113").Value = 2485.6667
$ws.Range("I113").Value = 2004
$ws.Range("J113").Value = 2582
$ws.Range("K113").Value = 2004
$ws.Range("L113").Value = 2582
$ws.Range("M113").Value = 166
$ws.Range("N113").Value = -6922

$ws.Range("H132").Value = 24128.71
$ws.Range("I132").Value = 1455.56
$ws.Range("J132").Value = 52470.15
$ws.Range("K132").Value = 4366.68
$ws.Range("L132").Value = 157410.45
$ws.Range("M132").Value = -1836.68
$ws.Range("N132").Value = -162470.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents() | Out-Null

$ws.Range("H113").Value = 398
$ws.Range("I113").Value = 398
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1194
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 976

$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("M114").ClearContents() | Out-Null

$ws.Range("H127").Value = 74500
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 74500
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 74500
$ws.Range("N127").Value = -84420
